$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: header "PSET_CO" -> "value"; drop the CSET_CN / value columns (L7, M7)
$ws.Range("K7").Value = "value"
$ws.Range("L7").Clear()
$ws.Range("M7").Clear()

# Row 9: add new D9 label first (matches original authoring order), drop the
# now-unused H9/J9/K9/M9 cells
$ws.Range("D9").Value = "*V2G"
$ws.Range("H9").Clear()
$ws.Range("J9").Clear()
$ws.Range("K9").Clear()
$ws.Range("M9").Clear()

# Row 8: add new D8 label, rework the formulas/values in G8..M8
$ws.Range("D8").Value = "*SIFTELC1"
$ws.Range("G8").Value = "UP"
$ws.Range("H8").Formula = '=IF($C$5="Ref","\I:","NCAP_BND")'
$ws.Range("J8").Formula = '=VLOOKUP(C5,$C$7:$D$9,2,FALSE)'
$ws.Range("K8").Value = 0
$ws.Range("L8").Clear()
$ws.Range("M8").Clear()

# Move the active selection to K9 (matches the saved view state)
$ws.Range("K9").Select()
